$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data rows 2-15 with new values ---
# Row 2
$ws.Range("A2").Value = "'Shofa Khairunnisa"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "'Anisa_s1"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "'s1"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value2 = 3
$ws.Range("E2").Value = "'556,912.00"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = "'57,795,497.00"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "'0.96"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value2 = 0
$ws.Range("I2").Value2 = 62
$ws.Range("J2").Value2 = 3
$ws.Range("K2").Value = "'1.40"
$ws.Range("K2").Style = "Normal"
$ws.Range("L2").Value = "'4.84"
$ws.Range("L2").Style = "Normal"

# Row 3
$ws.Range("A3").Value = "'Febri Fransiska"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "'Anisa_s1"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "'s1"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value2 = 12
$ws.Range("E3").Value = "'4,454,525.00"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = "'74,289,679.00"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = "'6.00"
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value2 = 227
$ws.Range("I3").Value2 = 62
$ws.Range("J3").Value2 = 3
$ws.Range("K3").Value = "'3.53"
$ws.Range("K3").Style = "Normal"
$ws.Range("L3").Value = "'4.84"
$ws.Range("L3").Style = "Normal"

# Row 4
$ws.Range("A4").Value = "'Ramesintia Sinaga"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = "'Anisa_s1"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "'s1"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value2 = 11
$ws.Range("E4").Value = "'6,373,174.00"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "'78,194,474.00"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = "'8.15"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value2 = 2.003
$ws.Range("I4").Value2 = 62
$ws.Range("J4").Value2 = 9
$ws.Range("K4").Value = "'10.31"
$ws.Range("K4").Style = "Normal"
$ws.Range("L4").Value = "'14.52"
$ws.Range("L4").Style = "Normal"

# Row 5
$ws.Range("A5").Value = "'Dwi Gusti Anggraini"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = "'Anisa_s1"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "'s1"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value2 = 17
$ws.Range("E5").Value = "'5,205,946.00"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "'57,062,552.00"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = "'9.12"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value2 = 775
$ws.Range("I5").Value2 = 62
$ws.Range("J5").Value2 = 6
$ws.Range("K5").Value = "'7.33"
$ws.Range("K5").Style = "Normal"
$ws.Range("L5").Value = "'9.68"
$ws.Range("L5").Style = "Normal"

# Row 6
$ws.Range("A6").Value = "'Okky Al Bana"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = "'Anisa_s1"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'s1"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value2 = 7
$ws.Range("E6").Value = "'3,819,840.00"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = "'65,002,470.00"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = "'5.88"
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Value2 = 0
$ws.Range("I6").Value2 = 62
$ws.Range("J6").Value2 = 3
$ws.Range("K6").Value = "'5.03"
$ws.Range("K6").Style = "Normal"
$ws.Range("L6").Value = "'4.84"
$ws.Range("L6").Style = "Normal"

# Row 7
$ws.Range("A7").Value = "'Dimas Kuat Anggowo"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = "'Anisa_s1"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'s1"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value2 = 14
$ws.Range("E7").Value = "'5,814,634.00"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "'64,783,410.00"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = "'8.98"
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value2 = 962
$ws.Range("I7").Value2 = 62
$ws.Range("J7").Value2 = 8
$ws.Range("K7").Value = "'5.60"
$ws.Range("K7").Style = "Normal"
$ws.Range("L7").Value = "'12.90"
$ws.Range("L7").Style = "Normal"

# Row 8
$ws.Range("A8").Value = "'Syaiful Munir"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = "'Anisa_s1"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'s1"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value2 = 10
$ws.Range("E8").Value = "'3,705,892.00"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = "'56,474,261.00"
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = "'6.56"
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").Value2 = 460
$ws.Range("I8").Value2 = 62
$ws.Range("J8").Value2 = 8
$ws.Range("K8").Value = "'7.78"
$ws.Range("K8").Style = "Normal"
$ws.Range("L8").Value = "'12.90"
$ws.Range("L8").Style = "Normal"

# Row 9
$ws.Range("A9").Value = "'Anisa Inraniwi"
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = "'Anisa_s1"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'s1"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value2 = 15
$ws.Range("E9").Value = "'7,197,709.00"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Value = "'64,703,721.00"
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Value = "'11.12"
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").Value2 = 928
$ws.Range("I9").Value2 = 65
$ws.Range("J9").Value2 = 4
$ws.Range("K9").Value = "'11.67"
$ws.Range("K9").Style = "Normal"
$ws.Range("L9").Value = "'6.15"
$ws.Range("L9").Style = "Normal"

# Row 10
$ws.Range("A10").Value = "'Fachrul Rozi"
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = "'Anisa_s1"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'s1"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value2 = 8
$ws.Range("E10").Value = "'3,926,689.00"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Value = "'61,578,870.00"
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").Value = "'6.38"
$ws.Range("G10").Style = "Normal"
$ws.Range("H10").Value2 = 875
$ws.Range("I10").Value2 = 62
$ws.Range("J10").Value2 = 1
$ws.Range("K10").Value = "'2.21"
$ws.Range("K10").Style = "Normal"
$ws.Range("L10").Value = "'1.61"
$ws.Range("L10").Style = "Normal"

# Row 11
$ws.Range("A11").Value = "'Gloriana Yesica"
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").Value = "'Anisa_s1"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'s1"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value2 = 3
$ws.Range("E11").Value = "'1,430,829.00"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = "'63,332,938.00"
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").Value = "'2.26"
$ws.Range("G11").Style = "Normal"
$ws.Range("H11").Value2 = 894
$ws.Range("I11").Value2 = 62
$ws.Range("J11").Value2 = 1
$ws.Range("K11").Value = "'0.68"
$ws.Range("K11").Style = "Normal"
$ws.Range("L11").Value = "'1.61"
$ws.Range("L11").Style = "Normal"

# Row 12
$ws.Range("A12").Value = "'Riska Rahmayanti"
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").Value = "'Anisa_s1"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'s1"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value2 = 7
$ws.Range("E12").Value = "'1,942,192.00"
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Value = "'65,635,601.00"
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").Value = "'2.96"
$ws.Range("G12").Style = "Normal"
$ws.Range("H12").Value2 = 460
$ws.Range("I12").Value2 = 63
$ws.Range("J12").Value2 = 5
$ws.Range("K12").Value = "'3.19"
$ws.Range("K12").Style = "Normal"
$ws.Range("L12").Value = "'7.94"
$ws.Range("L12").Style = "Normal"

# Row 13
$ws.Range("A13").Value = "'Raina Claresta Purwanika"
$ws.Range("A13").Style = "Normal"
$ws.Range("B13").Value = "'Anisa_s1"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'s1"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value2 = 9
$ws.Range("E13").Value = "'1,756,804.00"
$ws.Range("E13").Style = "Normal"
$ws.Range("F13").Value = "'57,104,183.00"
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").Value = "'3.08"
$ws.Range("G13").Style = "Normal"
$ws.Range("H13").Value2 = 1.479
$ws.Range("I13").Value2 = 63
$ws.Range("J13").Value2 = 3
$ws.Range("K13").Value = "'3.14"
$ws.Range("K13").Style = "Normal"
$ws.Range("L13").Value = "'4.76"
$ws.Range("L13").Style = "Normal"

# Row 14
$ws.Range("A14").Value = "'Ayu Lestari"
$ws.Range("A14").Style = "Normal"
$ws.Range("B14").Value = "'Anisa_s1"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'s1"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value2 = 9
$ws.Range("E14").Value = "'1,639,600.00"
$ws.Range("E14").Style = "Normal"
$ws.Range("F14").Value = "'62,147,309.00"
$ws.Range("F14").Style = "Normal"
$ws.Range("G14").Value = "'2.64"
$ws.Range("G14").Style = "Normal"
$ws.Range("H14").Value2 = 1.579
$ws.Range("I14").Value2 = 62
$ws.Range("J14").Value2 = 5
$ws.Range("K14").Value = "'3.53"
$ws.Range("K14").Style = "Normal"
$ws.Range("L14").Value = "'8.06"
$ws.Range("L14").Style = "Normal"

# Row 15
$ws.Range("A15").Value = "'Neneng Hikmatul"
$ws.Range("A15").Style = "Normal"
$ws.Range("B15").Value = "'Anisa_s1"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'s1"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value2 = 12
$ws.Range("E15").Value = "'4,313,670.00"
$ws.Range("E15").Style = "Normal"
$ws.Range("F15").Value = "'63,277,954.00"
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").Value = "'6.82"
$ws.Range("G15").Style = "Normal"
$ws.Range("H15").Value2 = 593
$ws.Range("I15").Value2 = 64
$ws.Range("J15").Value2 = 4
$ws.Range("K15").Value = "'8.33"
$ws.Range("K15").Style = "Normal"
$ws.Range("L15").Value = "'6.25"
$ws.Range("L15").Style = "Normal"

# --- Column widths (best-fit approximation) ---
$ws.Columns.Item(1).ColumnWidth = 22.666666666666668
$ws.Columns.Item(2).ColumnWidth = 7.833333333333333
$ws.Columns.Item(3).ColumnWidth = 4.833333333333334
$ws.Columns.Item(4).ColumnWidth = 21.333333333333336
$ws.Columns.Item(5).ColumnWidth = 18.5
$ws.Columns.Item(6).ColumnWidth = 15.166666666666666
$ws.Columns.Item(7).ColumnWidth = 24.0
$ws.Columns.Item(8).ColumnWidth = 8.833333333333332
$ws.Columns.Item(9).ColumnWidth = 15.0
$ws.Columns.Item(10).ColumnWidth = 26.333333333333336
$ws.Columns.Item(11).ColumnWidth = 27.0
$ws.Columns.Item(12).ColumnWidth = 25.166666666666668

# --- Selection ---
$ws.Range("B10").Select() | Out-Null
